$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.882.50'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '3.859.70'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '473.15'
$ws.Range('E5').Value = '  +10.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.11'
$ws.Range('E6').Value = '  +10.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +2.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.744'
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000315'
$ws.Range('E11').Value = '  -5.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.33'
$ws.Range('E12').Value = '  +4.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.41'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('D14').Value = '4.483.38'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.89'
$ws.Range('E15').Value = '  -3.18%  '
$ws.Range('D16').Value = '3.859.67'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.12'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('E19').Value = '  +3.67%  '
$ws.Range('D20').Value = '67.197.09'
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.25'
$ws.Range('E21').Value = '  +4.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.97'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('E23').Value = '  +6.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.21'
$ws.Range('E24').Value = '  +2.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.60'
$ws.Range('E25').Value = '  +9.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.92'
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.03'
$ws.Range('E27').Value = '  +7.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.95'
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.53'
$ws.Range('E29').Value = '  +2.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '727.80'
$ws.Range('E30').Value = '  +2.46%  '
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('E32').Value = '  +6.23%  '
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.87'
$ws.Range('E34').Value = '  +10.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.160'
$ws.Range('E35').Value = '  +7.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.47'
$ws.Range('E36').Value = '  +4.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.52'
$ws.Range('E38').Value = '  -5.67%  '
$ws.Range('E39').Value = '  +3.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.95'
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.345'
$ws.Range('E41').Value = '  +7.19%  '
$ws.Range('E42').Value = '  +3.11%  '
$ws.Range('D43').Value = '0.0₃0680'
$ws.Range('E43').Value = '  -6.08%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.57'
$ws.Range('E44').Value = '  +7.15%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.46'
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.17'
$ws.Range('E47').Value = '  +5.70%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.76'
$ws.Range('E48').Value = '  +4.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.19'
$ws.Range('E49').Value = '  -2.38%  '
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.44'
$ws.Range('E51').Value = '  +1.30%  '
